{"js": "// Rewrite the \"Aulas Online\" bullet paragraph: the sentence describing\n// scheduling/editing/deleting classes, filtering by state, searching by\n// name, and the UC-schedule-conflict confirmation is reworded/expanded.\nconst START_ANCHOR = \"possibilita o agendamento\";\nconst END_ANCHOR = \"esta a decorrer.\";\nconst NEW_TEXT = \"possibilita o agendamento de uma aula online, a altera\u00e7\u00e3o da informa\u00e7\u00e3o da mesma e a sua elimina\u00e7\u00e3o. \u00c9 poss\u00edvel tamb\u00e9m listar todas as aulas ou escolher apenas o estado que pretende (agendadas, a decorrer ou realizadas), ou ainda, procurar uma aula em espec\u00edfico a partir do campo da designa\u00e7\u00e3o, na qual, tamb\u00e9m \u00e9 apresentado, para al\u00e9m da aula, os dados relativamente ao n\u00famero de estudantes presentes na mesma, como, \u00e0 quantidade de acesso \u00e0 sua grava\u00e7\u00e3o. Ao se agendar uma aula, tem de se escolher qual a hora de in\u00edcio e a partir dessa hora \u00e9 calculada a hora de fim, isto, atrav\u00e9s de um c\u00e1lculo relacionado com o n\u00famero de horas da disciplina e pela hora de in\u00edcio inserida anteriormente. Existe uma confirma\u00e7\u00e3o que n\u00e3o possibilita agendar uma aula quando outra aula da mesma UC estivar agendada ou a decorrer. Outras das confirma\u00e7\u00f5es durante o processo de agendamento, s\u00e3o a confirma\u00e7\u00e3o do c\u00f3digo da UC, da designa\u00e7\u00e3o, e da quantidade restante de horas da disciplina. \";\n\nconst body = context.document.body;\n\nconst startResults = body.search(START_ANCHOR, { matchCase: true });\nconst endResults = body.search(END_ANCHOR, { matchCase: true });\nstartResults.load('items');\nendResults.load('items');\nawait context.sync();\n\nif (startResults.items.length === 0 || endResults.items.length === 0) {\n  throw new Error('Could not locate the paragraph text to update.');\n}\n\n// Build a single range spanning from the start anchor through the end\n// anchor (inclusive) and replace its contents in one shot so the\n// paragraph's other runs (the underlined \"Aulas Online\" title, the\n// leading bullet/tab) are left untouched.\nconst fullRange = startResults.items[0].expandTo(endResults.items[0]);\nfullRange.insertText(NEW_TEXT, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Rewrite the \"Aulas Online\" bullet paragraph: the sentence describing\n# scheduling/editing/deleting classes, filtering by state, searching by\n# name, and the UC-schedule-conflict confirmation is reworded/expanded.\n$d = $word.ActiveDocument\n\n$startAnchor = \"possibilita o agendamento\"\n$endAnchor   = \"esta a decorrer.\"\n\n# Locate the start of the paragraph's body text (right after \"Aulas Online, \")\n$startRange = $d.Content\n$null = $startRange.Find.Execute($startAnchor)\nif (-not $startRange.Find.Found) {\n    throw \"Could not locate start anchor text.\"\n}\n$startPos = $startRange.Start\n\n# Locate the end of the sentence about scheduling conflicts\n$endRange = $d.Content\n$null = $endRange.Find.Execute($endAnchor)\nif (-not $endRange.Find.Found) {\n    throw \"Could not locate end anchor text.\"\n}\n$endPos = $endRange.End\n\n# Build one contiguous range spanning the whole block and replace its text\n# in a single shot, leaving the underlined \"Aulas Online\" title and the\n# leading bullet/tab runs untouched.\n$target = $d.Range($startPos, $endPos)\n$target.Text = \"possibilita o agendamento de uma aula online, a altera\u00e7\u00e3o da informa\u00e7\u00e3o da mesma e a sua elimina\u00e7\u00e3o. \u00c9 poss\u00edvel tamb\u00e9m listar todas as aulas ou escolher apenas o estado que pretende (agendadas, a decorrer ou realizadas), ou ainda, procurar uma aula em espec\u00edfico a partir do campo da designa\u00e7\u00e3o, na qual, tamb\u00e9m \u00e9 apresentado, para al\u00e9m da aula, os dados relativamente ao n\u00famero de estudantes presentes na mesma, como, \u00e0 quantidade de acesso \u00e0 sua grava\u00e7\u00e3o. Ao se agendar uma aula, tem de se escolher qual a hora de in\u00edcio e a partir dessa hora \u00e9 calculada a hora de fim, isto, atrav\u00e9s de um c\u00e1lculo relacionado com o n\u00famero de horas da disciplina e pela hora de in\u00edcio inserida anteriormente. Existe uma confirma\u00e7\u00e3o que n\u00e3o possibilita agendar uma aula quando outra aula da mesma UC estivar agendada ou a decorrer. Outras das confirma\u00e7\u00f5es durante o processo de agendamento, s\u00e3o a confirma\u00e7\u00e3o do c\u00f3digo da UC, da designa\u00e7\u00e3o, e da quantidade restante de horas da disciplina. \"\n"}
